# Update the shared "Target cluster" label: "Inflammatory-Mac" -> "ECs"
# and refresh the recomputed TPM-derived statistics for the Cxcl13-Cxcr5 pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 9).Value = 0.8781048434890718
$ws.Cells.Item(2, 10).Value = 0.8781048434890719
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.5604183333333334
$ws.Cells.Item(2, 14).Value = 1.681255
$ws.Cells.Item(2, 15).Value = 0.3661514391314925
$ws.Cells.Item(2, 16).Value = 0.3661514391314925
$ws.Cells.Item(2, 17).Value = 1.977630367522222
$ws.Cells.Item(2, 18).Value = 17.7986733077
$ws.Cells.Item(2, 19).Value = 0.3215193521518576
$ws.Cells.Item(2, 20).Value = 0.3215193521518576

# Row 3
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 9).Value = 0.8781048434890718
$ws.Cells.Item(3, 10).Value = 0.8781048434890719
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.29479
$ws.Cells.Item(3, 14).Value = 0.88437
$ws.Cells.Item(3, 15).Value = 0.1926021622090153
$ws.Cells.Item(3, 16).Value = 0.1926021622090153
$ws.Cells.Item(3, 17).Value = 1.040268708866667
$ws.Cells.Item(3, 18).Value = 9.362418379799999
$ws.Cells.Item(3, 19).Value = 0.1691248915022042
$ws.Cells.Item(3, 20).Value = 0.1691248915022042

# Row 4
$ws.Cells.Item(4, 9).Value = 0.8781048434890718
$ws.Cells.Item(4, 10).Value = 0.8781048434890719
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.439932
$ws.Cells.Item(4, 14).Value = 1.319796
$ws.Cells.Item(4, 15).Value = 0.2874312372364616
$ws.Cells.Item(4, 16).Value = 0.2874312372364616
$ws.Cells.Item(4, 17).Value = 1.55245257176
$ws.Cells.Item(4, 18).Value = 13.97207314584
$ws.Cells.Item(4, 19).Value = 0.2523947615873934
$ws.Cells.Item(4, 20).Value = 0.2523947615873934

# Row 5
$ws.Cells.Item(5, 9).Value = 0.8781048434890718
$ws.Cells.Item(5, 10).Value = 0.8781048434890719
$ws.Cells.Item(5, 13).Value = 0.235424
$ws.Cells.Item(5, 14).Value = 0.706272
$ws.Cells.Item(5, 15).Value = 0.1538151614230307
$ws.Cells.Item(5, 16).Value = 0.1538151614230307
$ws.Cells.Item(5, 17).Value = 0.8307751976533333
$ws.Cells.Item(5, 18).Value = 7.47697677888
$ws.Cells.Item(5, 19).Value = 0.1350658382476167
$ws.Cells.Item(5, 20).Value = 0.1350658382476167

# Row 6
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.489861
$ws.Cells.Item(6, 8).Value = 1.469583
$ws.Cells.Item(6, 9).Value = 0.1218951565109281
$ws.Cells.Item(6, 10).Value = 0.1218951565109281
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.5604183333333334
$ws.Cells.Item(6, 14).Value = 1.681255
$ws.Cells.Item(6, 15).Value = 0.3661514391314925
$ws.Cells.Item(6, 16).Value = 0.3661514391314925
$ws.Cells.Item(6, 17).Value = 0.274527085185
$ws.Cells.Item(6, 18).Value = 2.470743766665
$ws.Cells.Item(6, 19).Value = 0.04463208697963484
$ws.Cells.Item(6, 20).Value = 0.04463208697963484

# Row 7
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.489861
$ws.Cells.Item(7, 8).Value = 1.469583
$ws.Cells.Item(7, 9).Value = 0.1218951565109281
$ws.Cells.Item(7, 10).Value = 0.1218951565109281
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.29479
$ws.Cells.Item(7, 14).Value = 0.88437
$ws.Cells.Item(7, 15).Value = 0.1926021622090153
$ws.Cells.Item(7, 16).Value = 0.1926021622090153
$ws.Cells.Item(7, 17).Value = 0.14440612419
$ws.Cells.Item(7, 18).Value = 1.29965511771
$ws.Cells.Item(7, 19).Value = 0.02347727070681108
$ws.Cells.Item(7, 20).Value = 0.02347727070681108

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.489861
$ws.Cells.Item(8, 8).Value = 1.469583
$ws.Cells.Item(8, 9).Value = 0.1218951565109281
$ws.Cells.Item(8, 10).Value = 0.1218951565109281
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.439932
$ws.Cells.Item(8, 14).Value = 1.319796
$ws.Cells.Item(8, 15).Value = 0.2874312372364616
$ws.Cells.Item(8, 16).Value = 0.2874312372364616
$ws.Cells.Item(8, 17).Value = 0.215505529452
$ws.Cells.Item(8, 18).Value = 1.939549765068
$ws.Cells.Item(8, 19).Value = 0.03503647564906819
$ws.Cells.Item(8, 20).Value = 0.03503647564906819

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.489861
$ws.Cells.Item(9, 8).Value = 1.469583
$ws.Cells.Item(9, 9).Value = 0.1218951565109281
$ws.Cells.Item(9, 10).Value = 0.1218951565109281
$ws.Cells.Item(9, 13).Value = 0.235424
$ws.Cells.Item(9, 14).Value = 0.706272
$ws.Cells.Item(9, 15).Value = 0.1538151614230307
$ws.Cells.Item(9, 16).Value = 0.1538151614230307
$ws.Cells.Item(9, 17).Value = 0.115325036064
$ws.Cells.Item(9, 18).Value = 1.037925324576
$ws.Cells.Item(9, 19).Value = 0.01874932317541399
$ws.Cells.Item(9, 20).Value = 0.01874932317541399
